$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$words = @(
    'töten',
    'scheitern',
    'prügeln',
    'sterben',
    'ehren',
    'seufzen',
    'wohnen',
    'fehlen',
    'bleiben',
    'hauen',
    'altern',
    'formen',
    'ärgern',
    'tollen',
    'fragen',
    'schrecken',
    'retten',
    'drehen',
    'spüren',
    'geben',
    'treiben',
    'platzen',
    'runden',
    'nennen',
    'wirken',
    'opfern',
    'drücken',
    'dringen',
    'fesseln',
    'feiern',
    'gelten',
    'enden',
    'pflanzen',
    'räumen',
    'ändern',
    'meistern',
    'schwingen',
    'drohen',
    'flüchten',
    'lockern',
    'wachsen',
    'fließen',
    'heilen',
    'feuern',
    'fischen',
    'planen',
    'klingen',
    'klappen',
    'münzen',
    'dauern',
    'filmen',
    'freuen',
    'siegen',
    'kennen',
    'trennen',
    'lügen',
    'kichern',
    'wehen',
    'hören',
    'schwimmen',
    'bluten',
    'stammen',
    'schenken',
    'liegen',
    'fangen',
    'stecken',
    'wellen',
    'wundern',
    'sperren',
    'graben',
    'doppeln',
    'schulden',
    'zögern',
    'folgen',
    'machen',
    'wiegen',
    'suchen',
    'grüßen',
    'sichern',
    'scheinen',
    'decken',
    'löschen',
    'hassen',
    'betteln',
    'stimmen',
    'greifen',
    'herzen',
    'schwören',
    'segeln',
    'schreiten',
    'stärken',
    'achten',
    'laufen',
    'sprengen',
    'tropfen',
    'erben',
    'sinken',
    'schmecken',
    'mauern',
    'leihen',
    'dienen',
    'rechnen',
    'irren',
    'klettern',
    'heulen',
    'schaden',
    'warten',
    'bellen',
    'weichen',
    'bergen',
    'gnaden',
    'stehlen',
    'zeigen',
    'messen',
    'plaudern',
    'rasen',
    'pfeifen',
    'biegen',
    'jagen',
    'wahren',
    'mögen',
    'pissen',
    'werden',
    'wenden',
    'boxen',
    'fühlen',
    'foltern',
    'sorgen',
    'brauchen',
    'führen',
    'schlucken',
    'fallen',
    'trauen',
    'stürzen',
    'knarren',
    'zielen',
    'lieben',
    'werfen',
    'bitten',
    'helfen',
    'flehen',
    'warnen',
    'tragen',
    'heben',
    'tauchen',
    'spielen',
    'fällen',
    'liefern',
    'drucken',
    'streichen',
    'zeugen',
    'wehtun',
    'saufen',
    'spinnen',
    'holen',
    'danken',
    'kosten',
    'narren',
    'streifen',
    'bauen',
    'reizen',
    'atmen',
    'kriegen',
    'lösen',
    'schwächen',
    'zahlen',
    'arten',
    'äußern',
    'quälen',
    'gründen',
    'schultern',
    'mühen',
    'wüten',
    'locken',
    'zünden',
    'lächeln',
    'kehren',
    'kümmern',
    'fahren',
    'husten',
    'lesen',
    'trösten',
    'ziehen',
    'wollen',
    'rufen',
    'rücken',
    'malen',
    'jubeln',
    'loben',
    'bilden',
    'weigern',
    'backen'
)
for ($i = 0; $i -lt $words.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}
